# Part1-UML-Review.pptx -- "adding team gradebook as dp2 makeup"
#
# The author filled in the "Today's Attendance password" textbox on the
# first slide: the blank-line placeholder ("__________") was replaced with
# the actual password text ("errorsasobjects").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Today's Attendance password" textbox by name (robust to any
# shape re-ordering) rather than by a hard-coded collection index.
$sh = $s.Shapes.Item("TextBox 1")

$tr = $sh.TextFrame.TextRange
$pass = $tr.Paragraphs(2)
$pass.Text = "errorsasobjects"
